$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their text (string) representation instead of
# being auto-converted to numbers/percentages by Excel when we assign values.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "59.091.15"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "2.639.74"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "527.12"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "144.47"
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("E9").Value = "  -4.25%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").Value = "0.336"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "3.114.01"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").Value = "59.050.94"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "20.98"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "2.663.32"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "341.76"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "4.45"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "10.55"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").Value = "6.34"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "65.46"
$ws.Range("E23").Value = "  +3.17%  "
$ws.Range("D24").Value = "0.418"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "7.24"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "0.0₃0799"
$ws.Range("E28").Value = "  -1.74%  "
$ws.Range("D29").Value = "6.47"
$ws.Range("E29").Value = "  -4.33%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("D32").Value = "18.95"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "150.11"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "4.21"
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("D37").Value = "0.874"
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").Value = "36.57"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  -6.01%  "
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").Value = "271.57"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").Value = "19.47"
$ws.Range("E45").Value = "  -3.27%  "
$ws.Range("D46").Value = "0.0539"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("D48").Value = "2.050.23"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").Value = "4.78"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0230"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "19.05"
$ws.Range("E51").Value = "  -0.27%  "

# Restore the original (default/Normal) style so no stray per-cell style
# attribute is introduced by the temporary text number format above.
$rng.Style = "Normal"

